# Generate Report for Handoff
# Regenerates the handoff report: a new handoff id/commit
# (3c147a42-66ec-4dcb-9945-7744755997f2) replaces the previous one
# (e2c76f03-beef-4c59-90eb-2b1b55f129fa), together with fresh handoff
# timestamps and the new xlf target-file names.

$wb = $excel.ActiveWorkbook

$oldId = "e2c76f03-beef-4c59-90eb-2b1b55f129fa"
$newId = "3c147a42-66ec-4dcb-9945-7744755997f2"

$oldZhXlf = "$oldId.d9038e2eea13e679c040f391382f7f83c3be626e.zh-cn.xlf"
$newZhXlf = "$newId.78825cd30660582252c00bc5f671e569db17818a.zh-cn.xlf"

$oldDeXlf = "$oldId.d9038e2eea13e679c040f391382f7f83c3be626e.de-de.xlf"
$newDeXlf = "$newId.78825cd30660582252c00bc5f671e569db17818a.de-de.xlf"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/c6ee848a901b26f924263c6eae85b19bb9473027/e2e/$oldId.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3dc95ff681b76ae03058344cfc0ccb09037abc27/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e893db80ffd2122ae9c57743239f4cd5ab87db5a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, "$newId.md")
$ws.Range("D2").Value = "2016-03-23 01:00:10"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, "$newId.md")
$ws.Hyperlinks.Add($ws.Range("D2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $newZhXlf)
$ws.Range("E2").Value = "2016-03-23 01:00:05"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [Type]::Missing, [Type]::Missing, "$newId.md")
$ws.Hyperlinks.Add($ws.Range("D2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $newDeXlf)
$ws.Range("E2").Value = "2016-03-23 01:00:10"
